$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.705.03"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.059.97"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.61"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.665"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.36"
$ws.Range("E8").Value = "  -6.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.64"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.368"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.956"
$ws.Range("E13").Value = "  +7.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.88"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").Value = "2.362.14"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("D17").Value = "2.060.88"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "36.619.71"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.48"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.30"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.84"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  -3.71%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("E26").Value = "  +5.48%  "
$ws.Range("E27").Value = "  -7.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.12"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.17"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("E31").Value = "  +7.40%  "
$ws.Range("E32").Value = "  -6.78%  "
$ws.Range("E33").Value = "  -4.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0598"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0851"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.21"
$ws.Range("E38").Value = "  -5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.09"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("E40").Value = "  -5.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("E41").Value = "  -5.95%  "
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("E43").Value = "  -4.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0917"
$ws.Range("E45").Value = "  -4.73%  "
$ws.Range("D46").Value = "1.417.67"
$ws.Range("E46").Value = "  +8.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.62"
$ws.Range("E47").Value = "  +11.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.04"
$ws.Range("E48").Value = "  -6.15%  "
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("E50").Value = "  -3.87%  "
$ws.Range("D51").Value = "2.247.09"
$ws.Range("E51").Value = "  +0.08%  "
